$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Period text (no date-like single token, so a
#     plain Value assignment is not re-interpreted as a date by Excel) ---
$wsMeta = $wb.Sheets.Item("Metadata")
$wsMeta.Range("J2").Value = "01-01-2024 to 31-03-2024"

# --- Transactions sheet: update each Transaction Date cell from
#     "DD-Mon-2024" to "DD-MM-2024". These values must stay plain TEXT
#     (as in the source file) rather than being auto-converted to Excel
#     date serials, so force Text number format before writing, then
#     restore the default "Normal" style afterwards so the cells keep the
#     same (unset) style index as before. ---
$wsTxn = $wb.Sheets.Item("Transactions")

$dates = $wsTxn.Range("A2:A13")
$dates.NumberFormat = "@"

$wsTxn.Range("A2").Value = "01-01-2024"
$wsTxn.Range("A3").Value = "03-01-2024"
$wsTxn.Range("A4").Value = "10-01-2024"
$wsTxn.Range("A5").Value = "15-01-2024"
$wsTxn.Range("A6").Value = "20-01-2024"
$wsTxn.Range("A7").Value = "28-01-2024"
$wsTxn.Range("A8").Value = "03-02-2024"
$wsTxn.Range("A9").Value = "12-02-2024"
$wsTxn.Range("A10").Value = "25-02-2024"
$wsTxn.Range("A11").Value = "01-03-2024"
$wsTxn.Range("A12").Value = "15-03-2024"
$wsTxn.Range("A13").Value = "31-03-2024"

$dates.Style = "Normal"
